# Re-upload edit: refresh the reported USD Amount figure and leave the
# workbook with the selection cursor on the updated cell's row (T3),
# matching the author's re-saved session state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edit: T2 224609 -> 225473
$ws.Range("T2").Value = 225473

# View edit: active cell / selection moves from T11 to T3
$ws.Activate() | Out-Null
$ws.Range("T3").Select() | Out-Null
